$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.225.75"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.886.43"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  -0.87%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.687"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.21%  "
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.06%  "
$ws.Range("E9").Value = "  -3.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0736"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0969"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").Value = "2.167.82"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.720"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.928.64"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "35.219.56"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.72%  "
$ws.Range("D20").Value = "0.0₃0818"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.46%  "
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("E25").Value = "  +6.84%  "
$ws.Range("E26").Value = "  -10.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("D31").Value = "4.128.40"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.83%  "
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0578"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.848"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -22.39%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0667"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("D45").Value = "1.292.16"
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0812"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.50%  "
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.46%  "
